# Update the OASH grant history table: the fiscal-year range referenced in the
# descriptive text moved from FY 2011-2016 to FY 2012-2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OASHgrantHistTable")

$ws.Activate()

$ws.Range("A3").Value = "This table shows the grant awards and award dollars OASH made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the OASH page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars OASH made for FY 2012-2016."

# Leave the active cell on A2, matching the saved selection state.
$ws.Range("A2").Select()
